# Updates cryptos list prices/volumes (and swaps Monero/EthereumClassic rows 27-28)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column to be treated as text so numeric-looking values
# (e.g. "315.47") are not auto-coerced to floating point numbers, matching
# the original inline-string cell content exactly.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '29.256.77'
$ws.Range("E2").Value = '  +2.40%  '
$ws.Range("D3").Value = '1.899.17'
$ws.Range("E3").Value = '  +0.66%  '
$ws.Range("D4").Value = '1.006'
$ws.Range("E4").Value = '  -0.92%  '
$ws.Range("D5").Value = '315.47'
$ws.Range("E5").Value = '  -0.36%  '
$ws.Range("D6").Value = '1.005'
$ws.Range("E6").Value = '  -0.71%  '
$ws.Range("D7").Value = '0.5138'
$ws.Range("E7").Value = '  +0.14%  '
$ws.Range("D8").Value = '0.3927'
$ws.Range("E8").Value = '  -1.29%  '
$ws.Range("D9").Value = '0.08435'
$ws.Range("E9").Value = '  -0.23%  '
$ws.Range("D10").Value = '42.54'
$ws.Range("E10").Value = '  +1.55%  '
$ws.Range("D11").Value = '1.116'
$ws.Range("E11").Value = '  +0.16%  '
$ws.Range("D12").Value = '6.252'
$ws.Range("E12").Value = '  -0.37%  '
$ws.Range("D13").Value = '1.898.50'
$ws.Range("E13").Value = '  +0.81%  '
$ws.Range("D14").Value = '20.71'
$ws.Range("E14").Value = '  +0.55%  '
$ws.Range("D15").Value = '7.322'
$ws.Range("E15").Value = '  +0.50%  '
$ws.Range("D16").Value = '1.007'
$ws.Range("E16").Value = '  -0.75%  '
$ws.Range("D17").Value = '93.43'
$ws.Range("E17").Value = '  +2.15%  '
$ws.Range("E18").Value = '  -0.36%  '
$ws.Range("D19").Value = '0.06740'
$ws.Range("E19").Value = '  -0.45%  '
$ws.Range("D20").Value = '17.86'
$ws.Range("E20").Value = '  +0.41%  '
$ws.Range("D21").Value = '1.003'
$ws.Range("E21").Value = '  -0.85%  '
$ws.Range("D22").Value = '6.016'
$ws.Range("E22").Value = '  +0.76%  '
$ws.Range("D23").Value = '29.240.92'
$ws.Range("E23").Value = '  +2.17%  '
$ws.Range("D24").Value = '11.15'
$ws.Range("E24").Value = '  -0.17%  '
$ws.Range("D25").Value = '2.219'
$ws.Range("E25").Value = '  -2.26%  '
$ws.Range("D26").Value = '2.118.19'
$ws.Range("E26").Value = '  +0.92%  '
$ws.Range("B27").Value = 'Monero'
$ws.Range("C27").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D27").Value = '159.33'
$ws.Range("E27").Value = '  -1.08%  '
$ws.Range("B28").Value = 'EthereumClassic'
$ws.Range("C28").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D28").Value = '20.91'
$ws.Range("E28").Value = '  +0.73%  '
$ws.Range("D29").Value = '2.448'
$ws.Range("E29").Value = '  +2.39%  '
$ws.Range("D30").Value = '128.09'
$ws.Range("E30").Value = '  +0.61%  '
$ws.Range("D31").Value = '1.058'
$ws.Range("E31").Value = '  +0.64%  '
$ws.Range("D32").Value = '0.1047'
$ws.Range("E32").Value = '  -1.04%  '
$ws.Range("D33").Value = '6.127'
$ws.Range("E33").Value = '  +5.67%  '
$ws.Range("D34").Value = '3.661'
$ws.Range("E34").Value = '  +1.21%  '
$ws.Range("D35").Value = '0.02479'
$ws.Range("E35").Value = '  +1.77%  '
$ws.Range("D36").Value = '0.06555'
$ws.Range("E36").Value = '  +0.81%  '
$ws.Range("D37").Value = '9.062'
$ws.Range("E37").Value = '  +1.05%  '
$ws.Range("D38").Value = '0.2193'
$ws.Range("E38").Value = '  +0.33%  '
$ws.Range("D39").Value = '1.231'
$ws.Range("E39").Value = '  +3.16%  '
$ws.Range("D40").Value = '5.133'
$ws.Range("E40").Value = '  +1.66%  '
$ws.Range("D41").Value = '0.6487'
$ws.Range("E41").Value = '  +0.65%  '
$ws.Range("D42").Value = '1.231'
$ws.Range("E42").Value = '  -2.76%  '
$ws.Range("D43").Value = '11.27'
$ws.Range("E43").Value = '  +0.32%  '
$ws.Range("D44").Value = '0.6056'
$ws.Range("E44").Value = '  -0.40%  '
$ws.Range("D45").Value = '13.14'
$ws.Range("E45").Value = '  +0.58%  '
$ws.Range("D46").Value = '3.680'
$ws.Range("E46").Value = '  -0.89%  '
$ws.Range("D47").Value = '2.047'
$ws.Range("E47").Value = '  +2.10%  '
$ws.Range("E48").Value = '  +1.47%  '
$ws.Range("D49").Value = '123.22'
$ws.Range("E49").Value = '  +0.46%  '
$ws.Range("D50").Value = '1.178'
$ws.Range("E50").Value = '  -2.63%  '
$ws.Range("D51").Value = '77.64'
$ws.Range("E51").Value = '  +0.61%  '

# Restore the original (default) cell style now that the text values are
# safely stored, so no stray number-format styling remains on the cells.
$ws.Range("D2:D51").Style = "Normal"
